# Add a new worksheet "Descriptiva_N12" (detailed per-run descriptive
# statistics, N=6 samples per group/metric) at the end of the workbook,
# right after the existing "Descriptiva_N2480" sheet.
#
# The new sheet has exactly the same shape/layout/header formatting as the
# existing sheet (same header row: Metrica, Grupo, N, Media, Mediana,
# Desv_Est, Min, Max, Q1, Q3), so we duplicate the existing sheet (which
# carries over its header style, page setup, column/row defaults, etc.
# faithfully) and then overwrite the data rows with the new numbers.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate the existing sheet and place the copy right after the last sheet.
$source.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Descriptiva_N12"

# New descriptive-statistics data (N=12 runs -> 6 per group) for each
# metric/group combination. Row 1 (headers) is left untouched - it was
# already copied over, style and all.
$rows = @(
    @("Instruction Coverage (%)", "Manual", 6, 18.25166666666667,  15.95,                12.50267717997496,  5.03,    35.56,     8.407500000000001,   27.3325),
    @("Instruction Coverage (%)", "IA",     6, 17.67525,           16.13,                11.36214431236463,  5.76,    33,        8.2225,              26.116125),
    @("Branch Coverage (%)",      "Manual", 6, 14.58333333333333,  11.875,               12.61612724518371,  2.5,     37.5,      6.5625,              17.1875),
    @("Branch Coverage (%)",      "IA",     6, 12.046875,          9.890625,             9.491644145971234,  2.5,     28.75,     5.8828125,           15),
    @("Mutation Score (%)",       "Manual", 6, 18.52,              15.28,                17.71322782555455,  0,       50,        7.6425,              22.9175),
    @("Mutation Score (%)",       "IA",     6, 14.757125,          15.106375,            11.61356661813028,  0,       33.33,     7.5556875,           18.7475),
    @("Time (seconds)",           "Manual", 6, 0.0822592261904762, 0.07102749999999999, 0.06882194001827036, 0.00932, 0.1884125, 0.03060089285714286, 0.12023875),
    @("Time (seconds)",           "IA",     6, 0.1937777083333333, 0.143878125,          0.1931202829042315, 0.005135, 0.433825, 0.0345515625,       0.36435)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    $rowIndex = $r + 2
    for ($col = 0; $col -lt $rowData.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $rowData[$col]
    }
}

# Leave the original sheet as the active/selected one, same as before the edit.
$wb.Worksheets.Item(1).Activate()
